{"js": "// Target paragraph texts, in final order (includes 4 new summary lines + 1 blank separator).\nconst targetParagraphs = [\n  \"LOTOFSAY: 0.32259444239638924\",\n  \"SETTHNGS: 0.2551938299003872,\",\n  \"TRAINOPS: 0.15561887504467253, \",\n  \"MYSKILLS: 0.11501904345844463,\",\n  \"WKDECIDE': -0.10881604027366734,\",\n  \"FAIREARN': -0.06298829981036164, \",\n  \"TRUSTMAN': 0.06070649684099467\",\n  \"LEARNNEW': 0.04813334995575874,\",\n  \"JOBSECOK': 0.04573300009304156,\",\n  \"WKFREEDM': -0.04163307291383008, \",\n  \"WKVSFAM': -0.04143096074950218, \",\n  \"WORKFAST': -0.03738125851839945,\",\n  \"KNOWWHAT': -0.02697069405608197, \",\n  \"WORKDIFF': -0.02535770010606593\",\n  \"SAFETYWK': -0.021078039254131733,\",\n  \"RESPECT': 0.01710195821765779,\",\n  \"HRSRELAX': 0.014934075248488\",\n  \"'HEALTH1': 0.008763002782717256\",\n  \"PROMTEOK': -0.005567600564531937, '\",\n  \"WKPRAISE': 0.0024888342128801043, '\",\n  \"\",\n  \"training data MSE 0.3237659428198579\",\n  \"test data MSE 0.29581819266384735\",\n  \"training data R-square 0.5689674198706031\",\n  \"test data R-square 0.5900191900772433\"\n];\n\nconst body = context.document.body;\n\n// Collapse the document down to its first paragraph; everything else gets rebuilt below.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst existing = body.paragraphs.items;\nfor (let i = existing.length - 1; i >= 1; i--) {\n  existing[i].delete();\n}\nawait context.sync();\n\n// Stamp the first target line into the lone remaining paragraph, then grow the rest of\n// the document by inserting a new paragraph after each one in turn.\nlet current = body.paragraphs.getFirst();\ncurrent.insertText(targetParagraphs[0], Word.InsertLocation.replace);\nawait context.sync();\n\nfor (let i = 1; i < targetParagraphs.length; i++) {\n  current = current.insertParagraph(targetParagraphs[i], Word.InsertLocation.after);\n  await context.sync();\n}\n\n// The \"_GoBack\" bookmark used to sit on the RESPECT line; re-home it on the HEALTH1 line,\n// which is where the author's cursor ended up after the edit.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst healthMatches = body.search(\"'HEALTH1': 0.008763002782717256\", { matchCase: true });\nhealthMatches.load(\"items\");\nawait context.sync();\nif (healthMatches.items.length > 0) {\n  healthMatches.items[0].insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the old _GoBack bookmark if present; we will re-add it at the correct spot.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Target paragraph texts, in final order (includes 4 new summary lines + 1 blank separator).\n$paragraphs = @(\n    \"LOTOFSAY: 0.32259444239638924\",\n    \"SETTHNGS: 0.2551938299003872,\",\n    \"TRAINOPS: 0.15561887504467253, \",\n    \"MYSKILLS: 0.11501904345844463,\",\n    \"WKDECIDE': -0.10881604027366734,\",\n    \"FAIREARN': -0.06298829981036164, \",\n    \"TRUSTMAN': 0.06070649684099467\",\n    \"LEARNNEW': 0.04813334995575874,\",\n    \"JOBSECOK': 0.04573300009304156,\",\n    \"WKFREEDM': -0.04163307291383008, \",\n    \"WKVSFAM': -0.04143096074950218, \",\n    \"WORKFAST': -0.03738125851839945,\",\n    \"KNOWWHAT': -0.02697069405608197, \",\n    \"WORKDIFF': -0.02535770010606593\",\n    \"SAFETYWK': -0.021078039254131733,\",\n    \"RESPECT': 0.01710195821765779,\",\n    \"HRSRELAX': 0.014934075248488\",\n    \"'HEALTH1': 0.008763002782717256\",\n    \"PROMTEOK': -0.005567600564531937, '\",\n    \"WKPRAISE': 0.0024888342128801043, '\",\n    \"\",\n    \"training data MSE 0.3237659428198579\",\n    \"test data MSE 0.29581819266384735\",\n    \"training data R-square 0.5689674198706031\",\n    \"test data R-square 0.5900191900772433\"\n)\n\n# Collapse the document down to a single paragraph, then stamp each target line in,\n# inserting a fresh paragraph mark after every line except the last.\nwhile ($d.Paragraphs.Count -gt 1) {\n    $d.Paragraphs(1).Range.Delete()\n}\n\n$count = $paragraphs.Count\nfor ($i = 0; $i -lt $count; $i++) {\n    $p = $d.Paragraphs($i + 1)\n    $p.Range.Text = $paragraphs[$i]\n    if ($i -lt ($count - 1)) {\n        $d.Paragraphs($i + 1).Range.InsertParagraphAfter()\n    }\n}\n\n# Re-attach the _GoBack bookmark to the HEALTH1 paragraph (its new home).\n$healthIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*HEALTH1*\") {\n        $healthIndex = $i\n        break\n    }\n}\nif ($healthIndex -gt 0) {\n    $healthRange = $d.Paragraphs($healthIndex).Range\n    $bookmarkRange = $d.Range($healthRange.Start, $healthRange.End - 1)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n}\n"}
